$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Interactive Presentation is an ideal instructional method for this course as it allows for dynamic engagement with the content related to storytelling and storyboarding using generative AI. This method enables the facilitator to present complex concepts such as AI-generated script ideation techniques and ethical considerations in a visually engaging manner. By incorporating multimedia elements and interactive components, learners can better grasp the narrative structures and creative storytelling principles essential for developing compelling script elements. This approach also encourages active participation, allowing learners to ask questions and engage in discussions, thereby enhancing their understanding of the material."

$ws.Range("B3").Value = "Brainstorming is a crucial instructional method for fostering creativity and innovation in the context of generative AI for storytelling and storyboarding. This method encourages learners to collaboratively generate ideas and explore various narrative components and prompt terms that can enhance AI-generated storyboards. By engaging in brainstorming sessions, learners can identify effective strategies for refining video scripts and developing unique story elements. This collaborative approach not only stimulates creative thinking but also helps learners to consider diverse perspectives, which is essential for creating inclusive and engaging narratives."

$ws.Range("B4").Value = "Demonstration is an effective instructional method for this course as it provides learners with a clear and concrete example of how to utilize generative AI tools for script and storyboard development. By visually showcasing the process of generating AI-driven story elements and refining video scripts, learners can better understand the practical application of AI techniques. This hands-on approach allows learners to see the results of their work in real-time, which can be motivating and rewarding. Demonstrations also help to break down complex processes into manageable steps, making it easier for learners to grasp the intricacies of AI-generated content creation."

$ws.Range("B5").Value = "Drill and Practice is an essential instructional method for reinforcing learners' understanding of generative AI tools and techniques. By providing opportunities for learners to repeatedly apply what they have learned, this method helps to solidify their skills in developing narrative structures and enhancing AI-generated storyboards. Through practice, learners can receive feedback on their performance, identify areas for improvement, and build confidence in their ability to create compelling and ethically sound AI-generated content. This method also prepares learners for future challenges by equipping them with the necessary skills and knowledge to succeed in the media industry."
